$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell A2 previously held the shared text "A001232" -> now reads "qqq"
$ws.Range("A2").Value = "qqq"

# Cell A3 previously held the number 4522342 -> now reads the text "www"
$ws.Range("A3").Value = "www"

# The sheet's active selection moved from D9 to E8
$ws.Range("E8").Select()
